$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts Agama..Alamat KTP one column to the right)
# and give it a header of "PTKP" for the new PTKP (non-taxable income) field
# used by the "transaction validation for role pegawai" work.
$ws.Columns("E:E").Insert()
$ws.Range("E1").Value = "PTKP"

# Recompute column widths for the (now 15-wide) header row, matching the
# widths Excel settled on after the column was inserted.
$ws.Columns("A:A").ColumnWidth = 13.166666666666666
$ws.Columns("B:B").ColumnWidth = 16.166666666666668
$ws.Columns("C:C").ColumnWidth = 16.666666666666668
$ws.Columns("D:E").ColumnWidth = 15.666666666666666
$ws.Columns("F:F").ColumnWidth = 12
$ws.Columns("G:G").ColumnWidth = 11.666666666666666
$ws.Columns("H:H").ColumnWidth = 17.333333333333332
$ws.Columns("I:I").ColumnWidth = 12.5
$ws.Columns("J:J").ColumnWidth = 16.166666666666668
$ws.Columns("K:K").ColumnWidth = 14.333333333333334
$ws.Columns("L:L").ColumnWidth = 5
$ws.Columns("M:M").ColumnWidth = 13.666666666666666
$ws.Columns("N:N").ColumnWidth = 6.5
$ws.Columns("O:O").ColumnWidth = 10.166666666666666

# The sheet view no longer needs to be scrolled to D1 - selecting G6 (which
# is already on-screen) both clears the stale topLeftCell and moves the
# active selection to where the editor left off.
$ws.Range("G6").Select()
